$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.357.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.833.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -7.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5183'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3237'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -8.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06737'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.76'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7675'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.85%  '

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.880.04'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07667'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.82'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.015'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.08'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007877'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.379.00'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.092.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.537'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.401'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -7.57%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.896'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.77'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.646'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.89'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '110.77'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.179'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.89%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.109'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.07%  '

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08717'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04826'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.127'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.847'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6802'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -8.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.097'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -6.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01782'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.195'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -9.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4893'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.95%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8986'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -8.14%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '111.29'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.01%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.707'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4177'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -9.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1256'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.34%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.093'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.32%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05876'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.29'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.87%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.416'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.38%  '
